$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 90.40000000000001
$ws.Range("I28").Value = 90.40000000000001
$ws.Range("K28").Value = 90.40000000000001
$ws.Range("M28").Value = 394.6
$ws.Range("H43").Value = 1542.8572
$ws.Range("J43").Value = 1520
$ws.Range("L43").Value = 1520
$ws.Range("N43").Value = -1658
$ws.Range("H51").Value = 7545.364
$ws.Range("J51").Value = 7374.875
$ws.Range("L51").Value = 7374.875
$ws.Range("N51").Value = -8342.875
$ws.Range("H58").Value = 285.16666
$ws.Range("I58").Value = 285.16666
$ws.Range("K58").Value = 855.4999799999999
$ws.Range("M58").Value = -705.4999799999999
$ws.Range("H116").Value = 36666
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H125").Value = 540.8889
$ws.Range("J125").Value = 550
$ws.Range("L125").Value = 4950
$ws.Range("N125").Value = -9870
$ws.Range("H136").Value = 69166.664
$ws.Range("J136").Value = 69166.664
$ws.Range("L136").Value = 69166.664
$ws.Range("N136").Value = -79366.664
$ws.Range("H137").Value = 30676.795
$ws.Range("I137").Value = 1236.7368
$ws.Range("K137").Value = 3710.2104
$ws.Range("M137").Value = -1160.2104
$ws.Range("H138").Value = 4975.159
$ws.Range("J138").Value = 4680.482
$ws.Range("L138").Value = 14041.446
$ws.Range("N138").Value = -24321.446
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12415.395
$ws.Range("I32").Value = 10567.265
$ws.Range("K32").Value = 10567.265
$ws.Range("M32").Value = -10280.265
$ws.Range("H45").Value = 2071.95
$ws.Range("I45").Value = 2212.375
$ws.Range("J45").Value = 1978.3334
$ws.Range("K45").Value = 2212.375
$ws.Range("L45").Value = 1978.3334
$ws.Range("M45").Value = -1835.375
$ws.Range("N45").Value = -2732.3334
$ws.Range("H135").Value = 37613.855
$ws.Range("J135").Value = 37613.855
$ws.Range("L135").Value = 37613.855
$ws.Range("N135").Value = -47753.855
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1499.0588
$ws.Range("I99").Value = 1298.8572
$ws.Range("K99").Value = 1298.8572
$ws.Range("M99").Value = 199.1428000000001
$ws.Range("H105").Value = 2388.6597
$ws.Range("I105").Value = 2366.2896
$ws.Range("J105").Value = 2483.111
$ws.Range("K105").Value = 2366.2896
$ws.Range("L105").Value = 2483.111
$ws.Range("M105").Value = -619.2896000000001
$ws.Range("N105").Value = -5977.111
$ws.Range("H107").Value = 1785.5714
$ws.Range("I107").Value = 1375
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 1375
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = 545
$ws.Range("N107").Value = -6173
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3811.1428
$ws.Range("I16").Value = 3811.1428
$ws.Range("K16").Value = 3811.1428
$ws.Range("M16").Value = -3524.1428
$ws.Range("H31").Value = 2516.5715
$ws.Range("I31").Value = 2405.5
$ws.Range("J31").Value = 2664.6667
$ws.Range("K31").Value = 2405.5
$ws.Range("L31").Value = 2664.6667
$ws.Range("M31").Value = -2110.5
$ws.Range("N31").Value = -3254.6667
$ws.Range("H34").Value = 2516.5715
$ws.Range("I34").Value = 2405.5
$ws.Range("J34").Value = 2664.6667
$ws.Range("K34").Value = 2405.5
$ws.Range("L34").Value = 2664.6667
$ws.Range("M34").Value = -2203.5
$ws.Range("N34").Value = -3068.6667
$ws.Range("H107").Value = 4600.3335
$ws.Range("J107").Value = 4273.75
$ws.Range("L107").Value = 4273.75
$ws.Range("N107").Value = -8113.75
$ws.Range("H113").Value = 3811.1428
$ws.Range("I113").Value = 3811.1428
$ws.Range("K113").Value = 3811.1428
$ws.Range("M113").Value = -1641.1428
$ws.Range("H141").Value = 152463
$ws.Range("J141").Value = 152463
$ws.Range("L141").Value = 152463
$ws.Range("N141").Value = -162823
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1477.6666
$ws.Range("J122").Value = 1949.5
$ws.Range("L122").Value = 17545.5
$ws.Range("N122").Value = -22445.5
$ws.Range("H132").Value = 1199.8334
$ws.Range("J132").Value = 1700
$ws.Range("L132").Value = 15300
$ws.Range("N132").Value = -20360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.14286
$ws.Range("I2").Value = 46.4
$ws.Range("K2").Value = 46.4
$ws.Range("M2").Value = 66.59999999999999
$ws.Range("H102").Value = 2548.4285
$ws.Range("I102").Value = 2347.8
$ws.Range("K102").Value = 2347.8
$ws.Range("M102").Value = -725.8000000000002
$ws.Range("H113").Value = 1606.1111
$ws.Range("I113").Value = 1347.3334
$ws.Range("J113").Value = 1864.8889
$ws.Range("K113").Value = 1347.3334
$ws.Range("L113").Value = 1864.8889
$ws.Range("M113").Value = 822.6666
$ws.Range("N113").Value = -6204.8889
$ws.Range("H132").Value = 1105258.6
$ws.Range("I132").Value = 1486280
$ws.Range("J132").Value = 4530.3335
$ws.Range("K132").Value = 4458840
$ws.Range("L132").Value = 13591.0005
$ws.Range("M132").Value = -4456310
$ws.Range("N132").Value = -18651.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14039.5
$ws.Range("J40").Value = 10769.077
$ws.Range("L40").Value = 10769.077
$ws.Range("N40").Value = -11041.077
$ws.Range("H46").Value = 3633
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H61").Value = 2324.7778
$ws.Range("I61").Value = 2109.4546
$ws.Range("J61").Value = 3272.2
$ws.Range("K61").Value = 2109.4546
$ws.Range("L61").Value = 3272.2
$ws.Range("M61").Value = -1907.4546
$ws.Range("N61").Value = -3676.2
$ws.Range("H93").Value = 20834362
$ws.Range("I93").Value = 1120.6923
$ws.Range("K93").Value = 1120.6923
$ws.Range("M93").Value = 127.3077000000001
$ws.Range("H113").Value = 2324.7778
$ws.Range("I113").Value = 2109.4546
$ws.Range("J113").Value = 3272.2
$ws.Range("K113").Value = 2109.4546
$ws.Range("L113").Value = 3272.2
$ws.Range("M113").Value = 60.54539999999997
$ws.Range("N113").Value = -7612.2
$ws.Range("H132").Value = 4326.591
$ws.Range("I132").Value = 3385.4666
$ws.Range("K132").Value = 10156.3998
$ws.Range("M132").Value = -7626.399800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 851.2
$ws.Range("I107").Value = 726.125
$ws.Range("J107").Value = 1351.5
$ws.Range("K107").Value = 2178.375
$ws.Range("L107").Value = 4054.5
$ws.Range("M107").Value = -258.375
$ws.Range("N107").Value = -7894.5
$ws.Range("H132").Value = 1422.3846
$ws.Range("I132").Value = 1044.6364
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 3133.9092
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -603.9092000000001
$ws.Range("N132").Value = -15560
